$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: Sunscreen ---
$ws.Range("A7").Value = "Sunscreen"
$ws.Range("B7").Value = 8
$ws.Range("C7").Formula = "=B7*TaxRate"
$ws.Range("D7").Formula = "=B7+C7"

# --- New row 8: Mentos ---
$ws.Range("A8").Value = "Mentos"
$ws.Range("B8").Value = 2
$ws.Range("C8").Formula = "=B8*TaxRate"
$ws.Range("D8").Formula = "=B8+C8"

# Give the new rows the same currency formatting as the rest of the table
# (reuses the existing "Currency" style instead of creating a new one)
$ws.Range("B7:D8").NumberFormat = $ws.Range("B6:D6").NumberFormat

# New running total in G6 using the (soon to be created) TaxPrice2 named range
$ws.Range("G6").Formula = "=SUM(TaxPrice2)"

# New named ranges / updated named range
$wb.Names.Add("NewProducts", "='Shopping Trip'!`$A`$7:`$A`$8")
$wb.Names.Add("TaxPrice2", "='Shopping Trip'!`$D`$2:`$D`$19")

# Leave the selection where the user ended up after adding the new total
$ws.Range("G6").Select()
